# Auto-generated: update Golem Profits market-data cells per scheduled runner diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 400
$ws.Range("K12").Value = 400
$ws.Range("M12").Value = -230
$ws.Range("H29").Value = 159.33333
$ws.Range("I29").Value = 64
$ws.Range("J29").Value = 350
$ws.Range("K29").Value = 192
$ws.Range("L29").Value = 1050
$ws.Range("M29").Value = 89
$ws.Range("N29").Value = -1612
$ws.Range("H39").Value = 24.642857
$ws.Range("I39").Value = 24.25
$ws.Range("J39").Value = 27
$ws.Range("K39").Value = 72.75
$ws.Range("L39").Value = 81
$ws.Range("M39").Value = 223.25
$ws.Range("N39").Value = -673
$ws.Range("H42").Value = 355.8
$ws.Range("I42").Value = 376.66666
$ws.Range("J42").Value = 324.5
$ws.Range("K42").Value = 1129.99998
$ws.Range("L42").Value = 973.5
$ws.Range("M42").Value = -899.9999800000001
$ws.Range("N42").Value = -1433.5
$ws.Range("H53").Value = 488.1111
$ws.Range("I53").Value = 465
$ws.Range("J53").Value = 569
$ws.Range("K53").Value = 465
$ws.Range("L53").Value = 569
$ws.Range("M53").Value = 172
$ws.Range("N53").Value = -1843
$ws.Range("H58").Value = 376.6
$ws.Range("I58").Value = 96.25
$ws.Range("J58").Value = 1498
$ws.Range("K58").Value = 288.75
$ws.Range("L58").Value = 4494
$ws.Range("M58").Value = -138.75
$ws.Range("N58").Value = -4794
$ws.Range("H86").Value = 981
$ws.Range("I86").Value = 981
$ws.Range("K86").Value = 981
$ws.Range("M86").Value = 142
$ws.Range("H89").Value = 981
$ws.Range("I89").Value = 981
$ws.Range("K89").Value = 4905
$ws.Range("M89").Value = 711
$ws.Range("H95").Value = 38000
$ws.Range("J95").Value = 38000
$ws.Range("L95").Value = 38000
$ws.Range("N95").Value = -43492
$ws.Range("H125").Value = 816.3333
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 724.5
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 6520.5
$ws.Range("M125").Value = -6540
$ws.Range("N125").Value = -11440.5
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H137").Value = 798.3333
$ws.Range("I137").Value = 798.3333
$ws.Range("K137").Value = 2394.9999
$ws.Range("M137").Value = 155.0001000000002
$ws.Range("H138").Value = 2584.25
$ws.Range("I138").Value = 2584.25
$ws.Range("K138").Value = 7752.75
$ws.Range("M138").Value = -2612.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1086.8
$ws.Range("J2").Value = 813
$ws.Range("L2").Value = 813
$ws.Range("N2").Value = -1039
$ws.Range("H102").Value = 252236
$ws.Range("J102").Value = 2978.3333
$ws.Range("L102").Value = 2978.3333
$ws.Range("N102").Value = -6222.3333
$ws.Range("H116").Value = 1086.8
$ws.Range("J116").Value = 813
$ws.Range("L116").Value = 813
$ws.Range("N116").Value = -5401

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1086.8
$ws.Range("J3").Value = 813
$ws.Range("L3").Value = 813
$ws.Range("N3").Value = -1041
$ws.Range("H86").Value = 1733.3334
$ws.Range("H89").Value = 1733.3334
$ws.Range("H107").Value = 2138.4443
$ws.Range("I107").Value = 2138.4443
$ws.Range("K107").Value = 2138.4443
$ws.Range("M107").Value = -218.4443000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1278
$ws.Range("I16").Value = 847.5
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 847.5
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -560.5
$ws.Range("N16").Value = -3574
$ws.Range("H88").Value = 43868.6
$ws.Range("J88").Value = 43868.6
$ws.Range("L88").Value = 43868.6
$ws.Range("N88").Value = -44680.6
$ws.Range("H91").Value = 43868.6
$ws.Range("J91").Value = 43868.6
$ws.Range("L91").Value = 43868.6
$ws.Range("N91").Value = -46676.6
$ws.Range("H112").Value = 75000
$ws.Range("J112").Value = 75000
$ws.Range("L112").Value = 75000
$ws.Range("N112").Value = -77954
$ws.Range("H113").Value = 1278
$ws.Range("I113").Value = 847.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 847.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1322.5
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 5219.769
$ws.Range("I132").Value = 3987.2
$ws.Range("K132").Value = 11961.6
$ws.Range("M132").Value = -9431.599999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 941.6
$ws.Range("J26").Value = 2241.5
$ws.Range("L26").Value = 6724.5
$ws.Range("N26").Value = -7300.5
$ws.Range("H33").Value = 414.58334
$ws.Range("I33").Value = 510.7143
$ws.Range("J33").Value = 280
$ws.Range("K33").Value = 3064.2858
$ws.Range("L33").Value = 1680
$ws.Range("M33").Value = -2781.2858
$ws.Range("N33").Value = -2246

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 74.62963000000001
$ws.Range("J2").Value = 106.666664
$ws.Range("L2").Value = 106.666664
$ws.Range("N2").Value = -332.666664
$ws.Range("H3").Value = 2857335.5
$ws.Range("I3").Value = 5000075
$ws.Range("J3").Value = 350
$ws.Range("K3").Value = 5000075
$ws.Range("L3").Value = 350
$ws.Range("M3").Value = -4999959
$ws.Range("N3").Value = -582
$ws.Range("H31").Value = 354.33334
$ws.Range("I31").Value = 354.33334
$ws.Range("K31").Value = 354.33334
$ws.Range("M31").Value = -62.33334000000002
$ws.Range("H37").Value = 354.33334
$ws.Range("I37").Value = 354.33334
$ws.Range("K37").Value = 354.33334
$ws.Range("M37").Value = -77.33334000000002
$ws.Range("H80").Value = 2242.5
$ws.Range("I80").Value = 2242.5
$ws.Range("K80").Value = 2242.5
$ws.Range("M80").Value = -1244.5
$ws.Range("H83").Value = 2242.5
$ws.Range("I83").Value = 2242.5
$ws.Range("K83").Value = 11212.5
$ws.Range("M83").Value = -6220.5
$ws.Range("H103").Value = 35000
$ws.Range("I103").Value = 35000
$ws.Range("K103").Value = 35000
$ws.Range("M103").Value = -33828
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070
$ws.Range("H122").Value = 4276.8823
$ws.Range("I122").Value = 3544.7
$ws.Range("J122").Value = 5322.857
$ws.Range("K122").Value = 10634.1
$ws.Range("L122").Value = 15968.571
$ws.Range("M122").Value = -8184.099999999999
$ws.Range("N122").Value = -20868.571
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -54900
$ws.Range("H126").Value = 8298.333000000001
$ws.Range("I126").Value = 8298.333000000001
$ws.Range("K126").Value = 24894.999
$ws.Range("M126").Value = -22424.999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1990
$ws.Range("J22").Value = 1990
$ws.Range("L22").Value = 1990
$ws.Range("N22").Value = -2580
$ws.Range("H27").Value = 1990
$ws.Range("J27").Value = 1990
$ws.Range("L27").Value = 1990
$ws.Range("N27").Value = -2204
$ws.Range("H68").Value = 3783
$ws.Range("I68").Value = 3124.75
$ws.Range("K68").Value = 3124.75
$ws.Range("M68").Value = -2375.75
$ws.Range("H71").Value = 3783
$ws.Range("I71").Value = 3124.75
$ws.Range("K71").Value = 15623.75
$ws.Range("M71").Value = -11879.75
$ws.Range("H132").Value = 2501500
$ws.Range("I132").Value = 5000000
$ws.Range("K132").Value = 15000000
$ws.Range("M132").Value = -14997470
$ws.Range("H136").Value = 1030299.4
$ws.Range("I136").Value = 1030299.4
$ws.Range("K136").Value = 3090898.2
$ws.Range("M136").Value = -3088348.2
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 4950
$ws.Range("J29").Value = 4950
$ws.Range("L29").Value = 4950
$ws.Range("N29").Value = -5530
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H124").Value = 33499.5
$ws.Range("J124").Value = 33499.5
$ws.Range("L124").Value = 33499.5
$ws.Range("N124").Value = -43319.5
$ws.Range("H132").Value = 1677.4
$ws.Range("I132").Value = 1349
$ws.Range("J132").Value = 1896.3334
$ws.Range("K132").Value = 4047
$ws.Range("L132").Value = 5689.0002
$ws.Range("M132").Value = -1517
$ws.Range("N132").Value = -10749.0002
$ws.Range("H136").Value = 1853.2142
$ws.Range("I136").Value = 1662.0834
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4986.2502
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2436.2502
$ws.Range("N136").Value = -14100
